# Remove trailing space in excel
# The cell that used to read "Result " (with a trailing space) should read "Result" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 2).Value = "Result"

# Move / update the current selection as recorded in the saved file.
$ws.Range("B12").Select()
